# Chart示例.xlsx edit: "Change StarHead Prefab, Add StarSound..."
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# taps (sheet2): tweak two B values
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("taps")
$ws2.Range("B60").Value = -0.7
$ws2.Range("B62").Value = -1

# ---------------------------------------------------------------
# holds (sheet3): LFunc/RFunc Sin -> Linear on two rows
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("holds")
$ws3.Range("I19").Value = "Linear"
$ws3.Range("J19").Value = "Linear"
$ws3.Range("I23").Value = "Linear"
$ws3.Range("J23").Value = "Linear"

# ---------------------------------------------------------------
# slides (sheet4): insert 4 rows at the top of the data block
# (old row 13 -> new row 17, i.e. everything from row 13 down
# shifts by +4) and append 5 new rows at the bottom.
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("slides")
$ws4.Rows.Item(13).Resize(4).Insert()

$ws4.Range("A13").Value = 9.686
$ws4.Range("B13").Value = 0.8
$ws4.Range("C13").Value = 0.6
$ws4.Range("D13").Value = 2

$ws4.Range("A14").Value = 9.736
$ws4.Range("B14").Value = 1.2
$ws4.Range("C14").Value = 0.4
$ws4.Range("D14").Value = 2

$ws4.Range("A15").Value = 9.686
$ws4.Range("B15").Value = -0.8
$ws4.Range("C15").Value = 0.6
$ws4.Range("D15").Value = 2

$ws4.Range("A16").Value = 9.736
$ws4.Range("B16").Value = -1.2
$ws4.Range("C16").Value = 0.4
$ws4.Range("D16").Value = 2

$ws4.Range("A13:A16").NumberFormat = "0.000"

$ws4.Range("A33").Value = 60.886000000000003 - 0.05
$ws4.Range("B33").Value = 0
$ws4.Range("C33").Value = 1.5
$ws4.Range("D33").Value = 2

$ws4.Range("A34").Value = 60.886
$ws4.Range("B34").Value = 0.8
$ws4.Range("C34").Value = 0.6
$ws4.Range("D34").Value = 2

$ws4.Range("A35").Value = (60.886000000000003 - 0.05) + 0.1
$ws4.Range("B35").Value = 1.2
$ws4.Range("C35").Value = 0.4
$ws4.Range("D35").Value = 2

$ws4.Range("A36").Value = 60.886
$ws4.Range("B36").Value = -0.8
$ws4.Range("C36").Value = 0.6
$ws4.Range("D36").Value = 2

$ws4.Range("A37").Value = (60.886000000000003 - 0.05) + 0.1
$ws4.Range("B37").Value = -1.2
$ws4.Range("C37").Value = 0.4
$ws4.Range("D37").Value = 2

$ws4.Range("A34:A37").NumberFormat = "0.000"

# ---------------------------------------------------------------
# flicks (sheet5): two Dir values
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("flicks")
$ws5.Range("D21").Value = 0.75
$ws5.Range("D25").Value = 0.25

# ---------------------------------------------------------------
# View state: update per-sheet selections (do the non-active
# sheets first so selecting on them doesn't "steal" the active
# tab), then finish on taps, which becomes the active sheet.
# ---------------------------------------------------------------
$ws3.Range("M23").Select()
$ws4.Range("L35").Select()
$ws5.Range("D25").Select()

$ws2.Activate()
$ws2.Range("G62").Select()
